$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column ("21-nov") before the
#     "01-oct." column (currently column DV), shifting everything from DV
#     onward one column to the right (DV:EZ -> DW:FA).
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("DV").Insert()

# New header cell for the inserted column.
$wsSpot.Range("DV1").Value = "21-nov"

# The inserted column has no data for this date yet -> "-" placeholder,
# matching every other not-yet-available day in the sheet.
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 126).Value = "-"
}

# --- Sheet "Gaz": append the next day's price.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A155").NumberFormat = "@"
$wsGaz.Range("A155").Value = "2025-11-19"
$wsGaz.Range("A155").Style = "Normal"
$wsGaz.Range("B155").Value = 29.925

# --- Sheet "CO2": append the next day's price.
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A155").NumberFormat = "@"
$wsCo2.Range("A155").Value = "2025-11-19"
$wsCo2.Range("A155").Style = "Normal"
$wsCo2.Range("B155").Value = 80.34
